$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.526.06"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "1.853.16"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.46"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4743"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2743"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06306"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.69"
$ws.Range("E10").Value = "  +10.17%  "
$ws.Range("D11").Value = "1.827.03"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07437"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6261"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "30.487.00"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "246.76"
$ws.Range("E17").Value = "  +9.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007324"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.932"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.909"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.121"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.52"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.98"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1020"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.352"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.018"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.827"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7007"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.701"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01896"
$ws.Range("E36").Value = "  +3.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("E38").Value = "  +3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8744"
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.45"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.539"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4046"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.191"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.68"
$ws.Range("E45").Value = "  +5.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1202"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.54"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.549"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05530"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.348"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3681"
$ws.Range("E51").Value = "  +0.97%  "
